$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '92.120.49'
$ws.Range('E2').Value = '  +0.54%  '

# Row 3
$ws.Range('D3').Value = '3.099.46'
$ws.Range('E3').Value = '  -1.28%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '233.62'
$ws.Range('E5').Value = '  -3.02%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '611.74'
$ws.Range('E6').Value = '  -1.25%  '

# Row 7
$ws.Range('E7').Value = '  -3.78%  '

# Row 8
$ws.Range('E8').Value = '  -0.75%  '

# Row 9
$ws.Range('E9').Value = '  -0.05%  '

# Row 10
$ws.Range('D10').Value = '3.094.79'
$ws.Range('E10').Value = '  -1.34%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.790'
$ws.Range('E11').Value = '  +5.39%  '

# Row 12
$ws.Range('E12').Value = '  -3.15%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000242'
$ws.Range('E13').Value = '  -4.94%  '

# Row 14
$ws.Range('D14').Value = '91.983.63'
$ws.Range('E14').Value = '  +0.78%  '

# Row 15
$ws.Range('B15').Value = 'Toncoin'
$ws.Range('C15').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.40'
$ws.Range('E15').Value = '  -3.91%  '

# Row 16
$ws.Range('B16').Value = 'Avalanche'
$ws.Range('C16').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '33.53'
$ws.Range('E16').Value = '  -4.43%  '

# Row 17
$ws.Range('D17').Value = '3.677.22'
$ws.Range('E17').Value = '  -1.04%  '

# Row 18
$ws.Range('D18').Value = '3.090.51'
$ws.Range('E18').Value = '  -3.33%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.78'
$ws.Range('E19').Value = '  +0.46%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.41'
$ws.Range('E20').Value = '  -3.61%  '

# Row 21
$ws.Range('E21').Value = '  -3.04%  '

# Row 22
$ws.Range('B22').Value = 'PEPE'
$ws.Range('C22').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.0000199'
$ws.Range('E22').Value = '  -1.64%  '

# Row 23
$ws.Range('B23').Value = 'BitcoinCash'
$ws.Range('C23').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '436.48'
$ws.Range('E23').Value = '  -4.48%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.10'
$ws.Range('E24').Value = '  -1.15%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.53'
$ws.Range('E25').Value = '  -6.58%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '85.10'
$ws.Range('E26').Value = '  -3.86%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.29'
$ws.Range('E27').Value = '  -4.76%  '

# Row 28
$ws.Range('D28').Value = '3.262.85'
$ws.Range('E28').Value = '  -1.65%  '

# Row 29
$ws.Range('E29').Value = '  +0.01%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.179'
$ws.Range('E30').Value = '  +7.78%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.232'
$ws.Range('E31').Value = '  +0.69%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.122'
$ws.Range('E32').Value = '  -19.76%  '

# Row 33
$ws.Range('B33').Value = 'InternetComputer(DFINITY)'
$ws.Range('C33').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '9.14'
$ws.Range('E33').Value = '  -2.30%  '

# Row 34
$ws.Range('B34').Value = 'Binance-PegBSC-USD'
$ws.Range('C34').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.999'
$ws.Range('E34').Value = '  +32.97%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '7.95'
$ws.Range('E35').Value = '  +6.82%  '

# Row 36
$ws.Range('E36').Value = '  -11.11%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '25.64'
$ws.Range('E37').Value = '  -2.87%  '

# Row 38
$ws.Range('E38').Value = '  -0.18%  '

# Row 39
$ws.Range('E39').Value = '  -5.53%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '23.84'
$ws.Range('E40').Value = '  +7.68%  '

# Row 41
$ws.Range('B41').Value = 'Fetch.AI'
$ws.Range('C41').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.27'
$ws.Range('E41').Value = '  -3.89%  '

# Row 42
$ws.Range('B42').Value = 'Bittensor'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '462.92'
$ws.Range('E42').Value = '  -6.02%  '

# Row 43
$ws.Range('B43').Value = 'PolygonEcosystemToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.430'
$ws.Range('E43').Value = '  -2.47%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.24'
$ws.Range('E44').Value = '  -4.52%  '

# Row 45
$ws.Range('E45').Value = '  +0.04%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '159.46'
$ws.Range('E46').Value = '  +1.76%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.679'
$ws.Range('E47').Value = '  -3.85%  '

# Row 48
$ws.Range('E48').Value = '  -5.37%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0324'
$ws.Range('E49').Value = '  -0.83%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '43.76'
$ws.Range('E50').Value = '  -0.69%  '

# Row 51
$ws.Range('E51').Value = '  -3.66%  '
